$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace cell A1 text content ("xyzz" -> "abcd")
$ws.Range("A1").Value = "abcd"

# Restore default selection to A1 (matches diff removing explicit A2 selection)
$ws.Range("A1").Select()
